$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    8 = 'dnasr281@gmail.com, System'
    9 = 'dnasr281@gmail.com, System'
    10 = 'dnasr281@gmail.com, System'
    12 = 'dnasr281@gmail.com, System'
    14 = 'dnasr281@gmail.com, System'
    15 = 'dnasr281@gmail.com, System'
    17 = 'dnasr281@gmail.com, System'
    18 = 'System, dnasr281@gmail.com'
    23 = 'System, dnasr281@gmail.com'
    34 = 'dnasr281@gmail.com, System'
    35 = 'dnasr281@gmail.com, System'
    36 = 'dnasr281@gmail.com, System'
    38 = 'dnasr281@gmail.com, System'
    40 = 'dnasr281@gmail.com, System'
    41 = 'dnasr281@gmail.com, System'
    43 = 'dnasr281@gmail.com, System'
    44 = 'System, dnasr281@gmail.com'
    49 = 'System, dnasr281@gmail.com'
    60 = 'dnasr281@gmail.com, System'
    61 = 'dnasr281@gmail.com, System'
    62 = 'dnasr281@gmail.com, System'
    64 = 'dnasr281@gmail.com, System'
    66 = 'dnasr281@gmail.com, System'
    67 = 'dnasr281@gmail.com, System'
    69 = 'dnasr281@gmail.com, System'
    70 = 'System, dnasr281@gmail.com'
    75 = 'System, dnasr281@gmail.com'
    86 = 'dnasr281@gmail.com, System'
    87 = 'dnasr281@gmail.com, System'
    88 = 'dnasr281@gmail.com, System'
    90 = 'dnasr281@gmail.com, System'
    92 = 'dnasr281@gmail.com, System'
    93 = 'dnasr281@gmail.com, System'
    95 = 'dnasr281@gmail.com, System'
    96 = 'System, dnasr281@gmail.com'
    101 = 'System, dnasr281@gmail.com'
    112 = 'dnasr281@gmail.com, System'
    113 = 'dnasr281@gmail.com, System'
    114 = 'dnasr281@gmail.com, System'
    116 = 'dnasr281@gmail.com, System'
    118 = 'dnasr281@gmail.com, System'
    119 = 'dnasr281@gmail.com, System'
    121 = 'dnasr281@gmail.com, System'
    122 = 'System, dnasr281@gmail.com'
    127 = 'System, dnasr281@gmail.com'
    138 = 'dnasr281@gmail.com, System'
    139 = 'dnasr281@gmail.com, System'
    140 = 'dnasr281@gmail.com, System'
    142 = 'dnasr281@gmail.com, System'
    144 = 'dnasr281@gmail.com, System'
    145 = 'dnasr281@gmail.com, System'
    147 = 'dnasr281@gmail.com, System'
    148 = 'System, dnasr281@gmail.com'
    153 = 'System, dnasr281@gmail.com'
    164 = 'System, dnasr281@gmail.com'
    167 = 'System, dnasr281@gmail.com'
    170 = 'System, dnasr281@gmail.com'
    174 = 'System, dnasr281@gmail.com'
    191 = 'System, dnasr281@gmail.com'
    194 = 'System, dnasr281@gmail.com'
    197 = 'System, dnasr281@gmail.com'
    201 = 'System, dnasr281@gmail.com'
    218 = 'System, dnasr281@gmail.com'
    221 = 'System, dnasr281@gmail.com'
    224 = 'System, dnasr281@gmail.com'
    228 = 'System, dnasr281@gmail.com'
    245 = 'System, dnasr281@gmail.com'
    248 = 'System, dnasr281@gmail.com'
    251 = 'System, dnasr281@gmail.com'
    255 = 'System, dnasr281@gmail.com'
    272 = 'System, dnasr281@gmail.com'
    275 = 'System, dnasr281@gmail.com'
    278 = 'System, dnasr281@gmail.com'
    282 = 'System, dnasr281@gmail.com'
    299 = 'System, dnasr281@gmail.com'
    302 = 'System, dnasr281@gmail.com'
    305 = 'System, dnasr281@gmail.com'
    309 = 'System, dnasr281@gmail.com'
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 7).Value = $updates[$row]
}
